# Sync automatico del tracker - updates completed-match results for
# rows 145-149 and appends 4 newly scraped fixtures (rows 160-163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a value into a cell while stopping Excel's automatic
# "looks like a date / percentage / number" conversion for strings that
# must stay literal text (e.g. "2025-09-18", "48.27%"). We temporarily
# force the cell to Text format, assign the literal value, then strip
# the formatting back off so the saved cell carries no explicit style.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

function Set-NumberValue {
    param($range, $value)
    $range.Value = $value
}

# ---------------------------------------------------------------------
# Results that came in for the matches played 2025-09-16 (rows 145-149)
# ---------------------------------------------------------------------

# Row 145 - PSV Eindhoven vs Union St. Gilloise -> Away Win (Fallo)
Set-TextValue $ws.Range("L145") "Completed"
Set-TextValue $ws.Range("M145") "Away Win"
Set-TextValue $ws.Range("N145") "Fallo"
Set-NumberValue $ws.Range("O145") -1.9
Set-NumberValue $ws.Range("P145") -100
Set-TextValue $ws.Range("Q145") "2025-09-17 04:26:45"

# Row 146 - Athletic Club vs Arsenal -> Away Win (Acierto)
Set-TextValue $ws.Range("L146") "Completed"
Set-TextValue $ws.Range("M146") "Away Win"
Set-TextValue $ws.Range("N146") "Acierto"
Set-NumberValue $ws.Range("O146") 1.7
Set-NumberValue $ws.Range("P146") 85
Set-TextValue $ws.Range("Q146") "2025-09-17 04:26:45"

# Row 147 - Juventus vs Borussia Dortmund -> Draw (Fallo)
Set-TextValue $ws.Range("L147") "Completed"
Set-TextValue $ws.Range("M147") "Draw"
Set-TextValue $ws.Range("N147") "Fallo"
Set-NumberValue $ws.Range("O147") -1.4
Set-NumberValue $ws.Range("P147") -100
Set-TextValue $ws.Range("Q147") "2025-09-17 04:26:45"

# Row 148 - Tottenham vs Villarreal -> Home Win (Acierto)
Set-TextValue $ws.Range("L148") "Completed"
Set-TextValue $ws.Range("M148") "Home Win"
Set-TextValue $ws.Range("N148") "Acierto"
Set-NumberValue $ws.Range("O148") 1.53
Set-NumberValue $ws.Range("P148") 85
Set-TextValue $ws.Range("Q148") "2025-09-17 04:26:45"

# Row 149 - Inter Miami vs Seattle Sounders -> Home Win (Acierto)
Set-TextValue $ws.Range("L149") "Completed"
Set-TextValue $ws.Range("M149") "Home Win"
Set-TextValue $ws.Range("N149") "Acierto"
Set-NumberValue $ws.Range("O149") 1.3
Set-NumberValue $ws.Range("P149") 100
Set-TextValue $ws.Range("Q149") "2025-09-17 04:26:45"

# ---------------------------------------------------------------------
# New predictions scraped for 2025-09-18 (appended as rows 160-163)
# ---------------------------------------------------------------------

$newRows = @(
    @{ Row=160; A="2025-09-18"; B="Major League Soccer";    C="Real Salt Lake";       D="Los Angeles FC";   E="Away Win"; F="48.27%"; G=2.25; H="7.51%";  I=0.4; J=0.006880097657915255; K=0.06880097657915255 },
    @{ Row=161; A="2025-09-18"; B="UEFA Champions League";  C="FC Copenhagen";        D="Bayer Leverkusen"; E="Away Win"; F="54.26%"; G=2.15; H="15.50%"; I=0.9; J=0.01449162217910136;  K=0.1449162217910136 },
    @{ Row=162; A="2025-09-18"; B="UEFA Champions League";  C="Eintracht Frankfurt";  D="Galatasaray";      E="Home Win"; F="54.05%"; G=2.1;  H="12.37%"; I=0.7; J=0.0122756764106551;   K=0.122756764106551 },
    @{ Row=163; A="2025-09-18"; B="UEFA Champions League";  C="Manchester City";      D="Napoli";           E="Home Win"; F="77.32%"; G=1.65; H="26.30%"; I=2.5; J=0.04242606754060104;  K=0.4242606754060104 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    Set-TextValue   $ws.Range("A$row") $r.A
    Set-TextValue   $ws.Range("B$row") $r.B
    Set-TextValue   $ws.Range("C$row") $r.C
    Set-TextValue   $ws.Range("D$row") $r.D
    Set-TextValue   $ws.Range("E$row") $r.E
    Set-TextValue   $ws.Range("F$row") $r.F
    Set-NumberValue $ws.Range("G$row") $r.G
    Set-TextValue   $ws.Range("H$row") $r.H
    Set-NumberValue $ws.Range("I$row") $r.I
    Set-NumberValue $ws.Range("J$row") $r.J
    Set-NumberValue $ws.Range("K$row") $r.K
    Set-TextValue   $ws.Range("L$row") "Pending"

    # M:Q stay blank (still "Pending"), but keep the cells materialised
    # -- present and empty -- mirroring every other pending row already
    # in the sheet (e.g. row 150..159) instead of leaving them absent.
    $blank = $ws.Range("M$row`:Q$row")
    $blank.NumberFormat = "@"
    $blank.Value = ""
    $blank.ClearFormats()
}
